# Generate Report for Handoff
# Adds two new source files (9631ccf0-...md and f05138b8-...md) to the
# localization status report: one new row per file on the "Overview"
# sheet, and one new row per file on each of the "zh-cn" / "de-de"
# per-locale detail sheets.

$wb = $excel.ActiveWorkbook

$linkColor = 15570276   # long(BGR) form of RGB(0x64,0x95,0xED) == FF6495ED

function Set-LinkFormat($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $linkColor
}

function Set-DateFormat($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewRows = @(
    @{ Row = 6; File = "9631ccf0-307f-4116-953f-8b31a329628a.md"; Status = "Ready for handoff"; Date = "2016-03-24 00:41:31" },
    @{ Row = 7; File = "f05138b8-5d19-4ce4-954b-9de98ee666c0.md"; Status = "Ready for handoff"; Date = "2016-03-24 00:41:31" }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $wsOverview.Range("A$row").Value = $r.File
    $wsOverview.Range("B$row").Value = $r.Status
    $wsOverview.Range("C$row").Value = $r.Status
    $wsOverview.Range("D$row").Value = $r.Date

    $wsOverview.Hyperlinks.Add($wsOverview.Range("A$row"), "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$($r.File)", "", "", $r.File)
    Set-LinkFormat $wsOverview.Range("A$row")
    Set-DateFormat $wsOverview.Range("D$row")
}

# ---------------------------------------------------------------------
# Per-locale sheets ("zh-cn" and "de-de"):
#   Source File Name | File Extension | Status | Latest Handoff File |
#   Latest Handoff Datetime | Latest Target File | Latest Handback File |
#   Latest Handback DateTime | Reference Tokens | Handoff Reason |
#   Dependency From | Error Detail
# ---------------------------------------------------------------------
$locales = @(
    @{
        Sheet = "zh-cn";
        Repo = "oltest.zh-cn";
        Rows = @(
            @{ Row = 6; File = "9631ccf0-307f-4116-953f-8b31a329628a.md"; Xlf = "9631ccf0-307f-4116-953f-8b31a329628a.6731550084e23c586360b5a7f352c43379e56c26.zh-cn.xlf"; HandoffDt = "2016-03-24 00:41:27" },
            @{ Row = 7; File = "f05138b8-5d19-4ce4-954b-9de98ee666c0.md"; Xlf = "f05138b8-5d19-4ce4-954b-9de98ee666c0.af946f76a5d932891c92feba719a00bb741806be.zh-cn.xlf"; HandoffDt = "2016-03-24 00:41:27" }
        )
    },
    @{
        Sheet = "de-de";
        Repo = "oltest.de-de";
        Rows = @(
            @{ Row = 6; File = "9631ccf0-307f-4116-953f-8b31a329628a.md"; Xlf = "9631ccf0-307f-4116-953f-8b31a329628a.6731550084e23c586360b5a7f352c43379e56c26.de-de.xlf"; HandoffDt = "2016-03-24 00:41:31" },
            @{ Row = 7; File = "f05138b8-5d19-4ce4-954b-9de98ee666c0.md"; Xlf = "f05138b8-5d19-4ce4-954b-9de98ee666c0.af946f76a5d932891c92feba719a00bb741806be.de-de.xlf"; HandoffDt = "2016-03-24 00:41:31" }
        )
    }
)

foreach ($locale in $locales) {
    $ws = $wb.Worksheets.Item($locale.Sheet)

    foreach ($r in $locale.Rows) {
        $row = $r.Row

        $ws.Range("A$row").Value = $r.File
        $ws.Range("B$row").Value = ".md"
        $ws.Range("C$row").Value = "Ready for handoff"
        $ws.Range("D$row").Value = $r.Xlf
        $ws.Range("E$row").Value = $r.HandoffDt
        $ws.Range("H$row").Value = "0001-01-01 00:00:00"
        $ws.Range("J$row").Value = "Include"

        $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/0000000000000000000000000000000000000000/e2e/$($r.File)"
        $xlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0000000000000000000000000000000000000000/ol-handoff/OpenLocalizationTestOrg/$($locale.Repo)/ci/ht/$($r.Xlf)"

        $ws.Hyperlinks.Add($ws.Range("A$row"), $mdUrl, "", "", $r.File)
        Set-LinkFormat $ws.Range("A$row")

        $ws.Hyperlinks.Add($ws.Range("D$row"), $xlfUrl, "", "", $r.Xlf)
        Set-LinkFormat $ws.Range("D$row")

        Set-DateFormat $ws.Range("E$row")
        Set-DateFormat $ws.Range("H$row")
    }
}

Write-Host "Report rows for handoff added."
